$d = $word.ActiveDocument

$d.Content.Find.Execute("354÷7=50, 4", $true, $false, $false, $false, $false, $true, 1, $false, "360÷3=120, 0", 2) | Out-Null
$d.Content.Find.Execute("416÷6=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "868÷5=173, 3", 2) | Out-Null
$d.Content.Find.Execute("887÷2=443, 1", $true, $false, $false, $false, $false, $true, 1, $false, "234÷9=26, 0", 2) | Out-Null
$d.Content.Find.Execute("444÷3=148, 0", $true, $false, $false, $false, $false, $true, 1, $false, "315÷2=157, 1", 2) | Out-Null
$d.Content.Find.Execute("898÷2=449, 0", $true, $false, $false, $false, $false, $true, 1, $false, "334÷4=83, 2", 2) | Out-Null
$d.Content.Find.Execute("861÷4=215, 1", $true, $false, $false, $false, $false, $true, 1, $false, "559÷3=186, 1", 2) | Out-Null
$d.Content.Find.Execute("539÷7=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "225÷2=112, 1", 2) | Out-Null
$d.Content.Find.Execute("994÷5=198, 4", $true, $false, $false, $false, $false, $true, 1, $false, "981÷4=245, 1", 2) | Out-Null
$d.Content.Find.Execute("901÷8=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "746÷5=149, 1", 2) | Out-Null
$d.Content.Find.Execute("760÷9=84, 4", $true, $false, $false, $false, $false, $true, 1, $false, "600÷6=100, 0", 2) | Out-Null
$d.Content.Find.Execute("359÷8=44, 7", $true, $false, $false, $false, $false, $true, 1, $false, "484÷3=161, 1", 2) | Out-Null
$d.Content.Find.Execute("461÷7=65, 6", $true, $false, $false, $false, $false, $true, 1, $false, "762÷3=254, 0", 2) | Out-Null
$d.Content.Find.Execute("192÷8=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "996÷3=332, 0", 2) | Out-Null
$d.Content.Find.Execute("541÷3=180, 1", $true, $false, $false, $false, $false, $true, 1, $false, "450÷3=150, 0", 2) | Out-Null
$d.Content.Find.Execute("966÷5=193, 1", $true, $false, $false, $false, $false, $true, 1, $false, "894÷9=99, 3", 2) | Out-Null
$d.Content.Find.Execute("620÷6=103, 2", $true, $false, $false, $false, $false, $true, 1, $false, "779÷7=111, 2", 2) | Out-Null
$d.Content.Find.Execute("330÷7=47, 1", $true, $false, $false, $false, $false, $true, 1, $false, "405÷7=57, 6", 2) | Out-Null
$d.Content.Find.Execute("513÷7=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "515÷7=73, 4", 2) | Out-Null
$d.Content.Find.Execute("237÷2=118, 1", $true, $false, $false, $false, $false, $true, 1, $false, "100÷7=14, 2", 2) | Out-Null
$d.Content.Find.Execute("223÷5=44, 3", $true, $false, $false, $false, $false, $true, 1, $false, "638÷2=319, 0", 2) | Out-Null
$d.Content.Find.Execute("740÷3=246, 2", $true, $false, $false, $false, $false, $true, 1, $false, "242÷7=34, 4", 2) | Out-Null
$d.Content.Find.Execute("401÷3=133, 2", $true, $false, $false, $false, $false, $true, 1, $false, "926÷6=154, 2", 2) | Out-Null
$d.Content.Find.Execute("118÷4=29, 2", $true, $false, $false, $false, $false, $true, 1, $false, "821÷4=205, 1", 2) | Out-Null
$d.Content.Find.Execute("943÷2=471, 1", $true, $false, $false, $false, $false, $true, 1, $false, "273÷2=136, 1", 2) | Out-Null
$d.Content.Find.Execute("696÷9=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "314÷9=34, 8", 2) | Out-Null
